# employees.xlsx — update the 4th data row (row 5) with Marina's real
# Telegram info (replacing the placeholder "@test3" row) and leave the
# newly-selected cell on A5, matching the author's last save state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tg_ID
$ws.Range("A5").Value = 836318110
# Tg_Username (new shared string "@IT_Svaha")
$ws.Range("B5").Value = "@IT_Svaha"
# BirthdayDate -> 2025-04-21 (serial date, keeps the existing custom
# datetime number format already applied to this cell)
$ws.Range("D5").Value = 45768
# NotificationTime -> 21:10:00 (fraction of a day, keeps existing time format)
$ws.Range("H5").Value = 0.88194444444444442

# Leave the selection on A5, as in the saved file.
$ws.Range("A5").Select() | Out-Null
